$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Countries table re-sorted: Kazajistan & Uzbekistan moved up in ranking,
#     pushing Azerbaiyan, Lituania(unchanged row), Armenia, Bosnia y Herzegovina,
#     Hong Kong and Republica de Macedonia(unchanged row) down by one position.
#     Update the country-name labels for the rows whose occupant changed.
$ws.Range("A70").Value = "Kazajistan"
$ws.Range("A71").Value = "Azerbaiyan"
$ws.Range("A73").Value = "Uzbekistan"
$ws.Range("A74").Value = "Armenia"
$ws.Range("A75").Value = "Bosnia y Herzegovina"
$ws.Range("A76").Value = "Hong Kong"

# --- Alemania (row 8)
$ws.Range("D8").Value = 68200
$ws.Range("E8").Value = 58678

# --- Brasil (row 17)
$ws.Range("D17").Value = 2979
$ws.Range("E17").Value = 19389

# --- India (row 25)
$ws.Range("B25").Value = 10541
$ws.Range("C25").Value = 88
$ws.Range("D25").Value = 1205
$ws.Range("E25").Value = 8978

# --- Pakistan (row 36)
$ws.Range("D36").Value = 1378
$ws.Range("E36").Value = 4233

# --- Hungria (row 63)
$ws.Range("B63").Value = 1512
$ws.Range("C63").Value = 54
$ws.Range("D63").Value = 122
$ws.Range("E63").Value = 1268
$ws.Range("G63").Value = 13
$ws.Range("H63").Value = 122

# --- Row 70 -> now Kazajistan
$ws.Range("B70").Value = 1179
$ws.Range("C70").Value = 88
$ws.Range("D70").Value = 138
$ws.Range("E70").Value = 1027
$ws.Range("F70").Value = 21
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 14

# --- Row 71 -> now Azerbaiyan
$ws.Range("B71").Value = 1148
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 289
$ws.Range("E71").Value = 847
$ws.Range("F71").Value = 25
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 12

# --- Row 73 -> now Uzbekistan
$ws.Range("B73").Value = 1054
$ws.Range("C73").Value = 56
$ws.Range("D73").Value = 85
$ws.Range("E73").Value = 965
$ws.Range("F73").Value = 8
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 4

# --- Row 74 -> now Armenia
$ws.Range("B74").Value = 1039
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 211
$ws.Range("E74").Value = 814
$ws.Range("F74").Value = 30
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 14

# --- Row 75 -> now Bosnia y Herzegovina
$ws.Range("B75").Value = 1037
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 206
$ws.Range("E75").Value = 792
$ws.Range("F75").Value = 4
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 39

# --- Row 76 -> now Hong Kong
$ws.Range("B76").Value = 1010
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 397
$ws.Range("E76").Value = 609
$ws.Range("F76").Value = 13
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 4

# --- Camboya (row 131)
$ws.Range("D131").Value = 91
$ws.Range("E131").Value = 31

# --- Curazao (row 187)
$ws.Range("D187").Value = 10
$ws.Range("E187").Value = 3

# --- Timestamp update
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 07:52"
